$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet tab to reflect new date
$ws.Name = "Through 2021-10-20"

# Row 9 (July 2021 columns T/U/V)
$ws.Range("U9").Value = 137
$ws.Range("V9").Value = 0.0805

# Row 12 (October, through 10-19 -> through 10-20)
$ws.Range("A12").Value = "October (through 10-20)"
$ws.Range("B12").Value = 2
$ws.Range("D12").Value = 0.1053
$ws.Range("F12").Value = 29
$ws.Range("G12").Value = 0.0938
$ws.Range("I12").Value = 28
$ws.Range("J12").Value = 0.2
$ws.Range("L12").Value = 45
$ws.Range("M12").Value = 0.0625
$ws.Range("N12").Value = 4
$ws.Range("O12").Value = 27
$ws.Range("P12").Value = 0.129
$ws.Range("R12").Value = 93
$ws.Range("U12").Value = 127

# Row 13 (Total)
$ws.Range("B13").Value = 32
$ws.Range("D13").Value = 0.1306
$ws.Range("F13").Value = 412
$ws.Range("G13").Value = 0.1063
$ws.Range("I13").Value = 605
$ws.Range("J13").Value = 0.0861
$ws.Range("L13").Value = 532
$ws.Range("M13").Value = 0.1074
$ws.Range("N13").Value = 47
$ws.Range("O13").Value = 406
$ws.Range("P13").Value = 0.1038
$ws.Range("R13").Value = 941
$ws.Range("S13").Value = 0.0533
$ws.Range("U13").Value = 1291
$ws.Range("V13").Value = 0.0604
